# Swap the match-data (columns B through AD) between specific pairs of
# adjacent rows. Column A (the sequential "id") is left untouched in
# each row; only the betting/match data that follows it is exchanged
# between the two rows of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$rowPairs = @(
    @(34,35),
    @(36,37),
    @(92,93),
    @(98,99),
    @(148,149),
    @(188,189),
    @(200,201),
    @(264,265),
    @(276,277),
    @(303,304),
    @(322,323)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $vals1 = @{}
    $vals2 = @{}

    foreach ($col in $cols) {
        $vals1[$col] = $ws.Range("$col$r1").Value2
        $vals2[$col] = $ws.Range("$col$r2").Value2
    }

    foreach ($col in $cols) {
        $ws.Range("$col$r1").Value2 = $vals2[$col]
        $ws.Range("$col$r2").Value2 = $vals1[$col]
    }
}
